$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game Clases")

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Apothecary"
$ws.Range("C17").Value = "chr"
$ws.Range("D17").Value = "focus"
$ws.Range("P17").Value = "Prophet"
$ws.Range("Q17").Value = "Merchant"
$ws.Range("R17").Value = 30
$ws.Range("S17").Value = 50
